# Atualização de bases das ligas, do dia: 26-02-2024 às 22:04
#
# This workbook has a single data sheet. Three pairs of rows had their
# match records swapped (everything except the running "id" in column A),
# and one row (121) had its result / closing-odds columns filled in for
# the first time (the match had not finished when the sheet was last saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 9 & 10: swap the two match records (column A "id" is untouched) ---
$range1 = $ws.Range("B9:AC9")
$range2 = $ws.Range("B10:AC10")
$tmp = $range1.Value()
$range1.Value = $range2.Value()
$range2.Value = $tmp

# --- Rows 87 & 88: swap the two match records ---
$range1 = $ws.Range("B87:AC87")
$range2 = $ws.Range("B88:AC88")
$tmp = $range1.Value()
$range1.Value = $range2.Value()
$range2.Value = $tmp

# --- Rows 99 & 100: swap the two match records ---
$range1 = $ws.Range("B99:AC99")
$range2 = $ws.Range("B100:AC100")
$tmp = $range1.Value()
$range1.Value = $range2.Value()
$range2.Value = $tmp

# --- Row 121: match result became available; fill in FTHG / FTAG / FTR
#     and the closing-line columns that were previously empty. ---
$ws.Range("H121").Value = 2
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = "H"

$ws.Range("W121").Value = 1.15
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = -1
$ws.Range("Z121").Value = 0.8999999999999999
$ws.Range("AA121").Value = -1
$ws.Range("AB121").Value = 0
$ws.Range("AC121").Value = -0
